# Edit script: merges the split word "t"/"he" (which were separated by the
# _GoBack bookmark) into a single run, then presses Enter at the end of the
# paragraph to start a new (empty) list item, moving the _GoBack bookmark
# there -- matching the author finishing step four of the problem solving
# list and beginning a new (not yet written) step.

$d = $word.ActiveDocument

# --- Step 1: merge the two runs that the _GoBack bookmark split in two ---
# Locate the text that spans from after the comma run through to the end
# of the paragraph (this text is unaffected in content, only run/bookmark
# structure changes).
$rng = $d.Content
$needle = " as there are no variables because she always counts the same exact way every time she is counting. This means that using this equation will give you the same results as if you were counting on your fingers the exact same way. "
$found = $rng.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$segStart = $rng.Start
$segEnd = $rng.End

# Remove the _GoBack bookmark (it currently sits between the two runs,
# contributing no visible text).
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# Re-insert the same text as a single clean run, replacing the old
# (now bookmark-free) two-run span.
$segRange = $d.Range($segStart, $segEnd)
$segText = $segRange.Text
$segRange.Delete()
$insertPoint = $d.Range($segStart, $segStart)
$insertPoint.InsertAfter($segText)

# --- Step 2: split the paragraph at its very end, creating a new, empty
# list item (same ListParagraph / numPr formatting) ---
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$pEnd = $lastPara.Range.End
$splitPoint = $d.Range($pEnd - 1, $pEnd - 1)
$splitPoint.InsertParagraphAfter()

# --- Step 3: put the _GoBack bookmark back, now at the start of the new,
# empty paragraph. A zero-length range placed directly at the start of a
# completely empty paragraph cannot be bookmarked directly in this host,
# so temporarily insert a placeholder character, anchor the bookmark
# there, then remove the placeholder again (the bookmark, being
# zero-length and anchored before the placeholder, survives). ---
$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newStart = $newPara.Range.Start
$placeholderPoint = $d.Range($newStart, $newStart)
$placeholderPoint.InsertAfter("X")

$newParaAfter = $d.Paragraphs($d.Paragraphs.Count)
$bmRange = $d.Range($newParaAfter.Range.Start, $newParaAfter.Range.Start)
$d.Bookmarks.Add("_GoBack", $bmRange)

$placeholderRange = $d.Range($newParaAfter.Range.Start, $newParaAfter.Range.Start + 1)
$placeholderRange.Delete()
